$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.026.39"
$ws.Range("E2").Value = "  +1.13%  "
$ws.Range("D3").Value = "3.380.48"
$ws.Range("E3").Value = "  +0.97%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.59%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "178.21"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.51%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.591"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.197"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +7.60%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.584"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.51%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "48.11"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.28%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000281"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.88%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "681.63"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.77%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "3.932.35"
$ws.Range("E14").Value = "  +1.11%  "
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.55"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.61%  "
$ws.Range("D16").Value = "69.311.45"
$ws.Range("E16").Value = "  +1.46%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.385.62"
$ws.Range("E17").Value = "  +0.75%  "
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.120"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.76%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.60"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.23"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.44%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.905"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.29%  "
$ws.Range("B22").Value = "Toncoin"
$ws.Range("C22").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.35"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.85%  "
$ws.Range("B23").Value = "InternetComputer(DFINITY)"
$ws.Range("C23").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.04"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.63%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "100.91"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.77%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.87"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.75%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.69"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.15%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.67"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.77%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "33.33"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.03%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.69"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.24%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.91"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.22%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.79"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +15.33%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.53%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "550.12"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.07%  "
$ws.Range("E34").Value = "  -0.15%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "57.71"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.30%  "
$ws.Range("E36").Value = "  +0.07%  "
$ws.Range("D37").Value = "3.593.60"
$ws.Range("E37").Value = "  -3.14%  "
$ws.Range("E38").Value = "  +2.69%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "35.19"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.44%  "
$ws.Range("D40").Value = "0.0₃0737"
$ws.Range("E40").Value = "  +10.01%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.29"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.35%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.70"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.73%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.37"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.17%  "
$ws.Range("E44").Value = "  +3.20%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.334"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.41%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.66"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.71%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.128"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.13%  "
$ws.Range("E48").Value = "  +3.60%  "
$ws.Range("E49").Value = "  -0.19%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "129.90"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.14%  "
$ws.Range("E51").Value = "  +1.62%  "
